# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The value is a date serial number that is bumped by one day
# (2026-02-22 -> 2026-02-23), i.e. 46075 -> 46076.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
